# "modify add member scenario"
# Adds a new member row to the MemberData sheet: A2 = "First9668".
# This mirrors the OOXML diff, which adds a new <row r="2"> to
# xl/worksheets/sheet3.xml whose A2 cell is a shared-string reference to
# the newly introduced "First9668" string (and bumps the sheet's
# <dimension> from A1 to A1:A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MemberData")
$ws.Range("A2").Value = "First9668"
